$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.827.75"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").Value = "2.221.00"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'228.60"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D7").Value = "'65.01"
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.408"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").Value = "'0.0880"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "2.547.35"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "'16.10"
$ws.Range("E13").Value = "  +0.76%  "
$ws.Range("D14").Value = "'22.39"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "'5.64"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "2.217.59"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "40.689.44"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("D19").Value = "'74.01"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +6.13%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'255.48"
$ws.Range("E22").Value = "  +10.42%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  -10.19%  "
$ws.Range("D26").Value = "'9.73"
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("D27").Value = "'173.08"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'0.146"
$ws.Range("E28").Value = "  +4.82%  "
$ws.Range("D29").Value = "'20.36"
$ws.Range("E29").Value = "  +2.16%  "
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "'2.86"
$ws.Range("E31").Value = "  +5.94%  "
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.83"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("B35").Value = "THORChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D35").Value = "'7.20"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "'3.82"
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("D38").Value = "'2.48"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'4.88"
$ws.Range("E40").Value = "  +14.85%  "
$ws.Range("D41").Value = "'0.0236"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("E42").Value = "  +10.15%  "
$ws.Range("D43").Value = "'102.07"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.24"
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("B45").Value = "TerraClassic"
$ws.Range("C45").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D45").Value = "'0.000220"
$ws.Range("E45").Value = "  +47.18%  "
$ws.Range("D46").Value = "'17.50"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "1.515.06"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").Value = "'0.0943"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.12"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.85"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("D51").Value = "'51.64"
$ws.Range("E51").Value = "  +11.58%  "
